# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.623.36"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "2.436.56"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'566.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("D6").Value = "'145.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("E9").Value = "  +2.48%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("E12").Value = "  +2.55%  "

$ws.Range("D13").Value = "'26.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.44%  "

$ws.Range("E14").Value = "  +6.11%  "

$ws.Range("D15").Value = "2.887.59"
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("D16").Value = "62.434.00"
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("D17").Value = "2.437.53"
$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").Value = "'324.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("E23").Value = "  +8.04%  "

$ws.Range("D24").Value = "'67.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.96%  "

$ws.Range("D25").Value = "'8.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D26").Value = "'580.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.25%  "

$ws.Range("E27").Value = "  +9.56%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("E30").Value = "  +4.28%  "

$ws.Range("E31").Value = "  +5.24%  "

$ws.Range("D32").Value = "'0.145"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("E33").Value = "  +0.44%  "

$ws.Range("E34").Value = "  +1.22%  "

$ws.Range("D35").Value = "'4.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.82%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("E37").Value = "  +1.42%  "

$ws.Range("D38").Value = "'18.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").Value = "'148.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.41%  "

$ws.Range("E41").Value = "  +2.82%  "

$ws.Range("D43").Value = "'2.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.61%  "

$ws.Range("D44").Value = "'148.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("E45").Value = "  +2.66%  "

$ws.Range("D46").Value = "'0.0534"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "

$ws.Range("D47").Value = "'20.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.49%  "

$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("E49").Value = "  +3.61%  "

$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("E51").Value = "  +5.15%  "
